$d = $word.ActiveDocument

# --- Register the built-in "Table Grid" table style in styles.xml -------
# (A side-effect of the original authoring session touching the table
#  gallery; materialise the style definition via a throw-away table so
#  it lands in word/styles.xml without altering the footer table below,
#  which keeps the default "TableNormal" style.)
$tmpRange = $d.Range(0, 0)
$tmpTable = $d.Tables.Add($tmpRange, 1, 1)
$tmpTable.Style = "Table Grid"
$tmpTable.Delete()

# --- Force an explicit page orientation (adds w:orient="portrait") ------
$d.PageSetup.Orientation = 0

# --- Split the header's tab + text runs into separate runs --------------
$hdr = $d.Sections.First.Headers.Item(1)
$hdrXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Header"/></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>Gosar</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:tab/></w:r><w:r><w:t>Use of Military Force in Syria</w:t></w:r><w:r><w:tab/></w:r><w:r><w:t>Sep 10, 2013</w:t></w:r></w:p>'
$hdr.Range.InsertXML($hdrXml)

# --- Add a three-column footer (blank "Three Columns" gallery layout) ---
$ftr = $d.Sections.First.Footers.Item(1)
$ftrXml = '<w:tbl xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:tblPr><w:tblStyle w:val="TableNormal"/><w:bidiVisual w:val="0"/><w:tblW w:w="0" w:type="auto"/><w:tblLayout w:type="fixed"/><w:tblLook w:val="04A0" w:firstRow="1" w:lastRow="0" w:firstColumn="1" w:lastColumn="0" w:noHBand="0" w:noVBand="1"/></w:tblPr><w:tblGrid><w:gridCol w:w="3120"/><w:gridCol w:w="3120"/><w:gridCol w:w="3120"/></w:tblGrid><w:tr><w:tc><w:tcPr><w:tcW w:w="3120" w:type="dxa"/><w:tcMar/></w:tcPr><w:p><w:pPr><w:pStyle w:val="Header"/><w:bidi w:val="0"/><w:ind w:left="-115"/><w:jc w:val="left"/></w:pPr></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3120" w:type="dxa"/><w:tcMar/></w:tcPr><w:p><w:pPr><w:pStyle w:val="Header"/><w:bidi w:val="0"/><w:jc w:val="center"/></w:pPr></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3120" w:type="dxa"/><w:tcMar/></w:tcPr><w:p><w:pPr><w:pStyle w:val="Header"/><w:bidi w:val="0"/><w:ind w:right="-115"/><w:jc w:val="right"/></w:pPr></w:p></w:tc></w:tr></w:tbl><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Footer"/><w:bidi w:val="0"/></w:pPr></w:p>'
$ftr.Range.InsertXML($ftrXml)

Write-Output "edit complete"
